$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (rows 2-6) date value from 45170 to 45174,
# preserving existing cell style/number format.
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45174
}
